# Issue Tracker.xlsx - Bugfix: add missing issue row for screen-size script bug
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new issue as row 8 (columns A = Issue Summary, D = Detail)
$ws.Range("A8").Value = "Script doesn't work for screen size 1102 by 677"
$ws.Range("D8").Value = "Need to make script work for that screen size"

# Move the selection to the newly added Detail cell, matching the author's
# on-screen position when they finished editing.
$ws.Range("D8").Select()
